$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Contact" (sheet1): restructure the columns, drop two obsolete
# contact rows (Jammer / Sample) and refresh the remaining three rows with
# the new FullName / ContactType / PinCode data.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Contact")

# Old layout:  A..C | D=ContactType | E=Email | F=Phone | G..Q | R=ProjectNotifyRole | S=Mentor
# New layout:  A..C | D=Email | E=Phone | F=FullName(new) | G=ContactType | H..Q | R=LineOfBusiness | S=PinCode

# Drop the old "ContactType" column (D) - Email/Phone shift left into D/E.
$ws1.Columns.Item(4).Delete()

# Insert two fresh columns at F/G for the new "FullName" and relocated
# "ContactType" fields.
$ws1.Columns.Item(6).Insert()
$ws1.Columns.Item(6).Insert()

# Drop the old "ProjectNotifyRole" column, which is now column S (19) -
# the remaining "Mentor" column shifts left to become the new column S,
# and gets overwritten below with "PinCode" data.
$ws1.Columns.Item(19).Delete()

# Remove the two contact rows that were dropped from the sheet (Jammer,
# Sample) - rows 5 and 6.
$ws1.Range("A5:A6").EntireRow.Delete()

# The freshly inserted F/G columns inherited the neighbouring cells'
# styling; reset them back to the default (unstyled) look before filling
# them in.
$ws1.Range("F2:G4").ClearFormats()

# Header row.
$ws1.Range("F1").Value = "FullName"
$ws1.Range("G1").Value = "ContactType"
$ws1.Range("S1").Value = "PinCode"

# Row 2 - Summer.
$ws1.Range("F2").Value = "Test CK Summer"
$ws1.Range("G2").Value = "External Contact"
$ws1.Range("Q2").Value = "CF"
$ws1.Range("R2").Value = "CF"
$ws1.Range("S2").Value = 92001

# Row 3 - Winter (also picks up the Country/Status/Office/... values that
# were previously blank on this row).
$ws1.Range("F3").Value = "Testing PS Winter"
$ws1.Range("G3").Value = "Distribution Lists"
$ws1.Range("L3").Value = "United States"
$ws1.Range("M3").Value = "Active"
$ws1.Range("N3").Value = "LA"
$ws1.Range("O3").Value = "DC"
$ws1.Range("P3").Value = "Associate"
$ws1.Range("Q3").Value = "CF"
$ws1.Range("R3").Value = "CF"
$ws1.Range("S3").Value = 92001

# Row 4 - Bingo (same treatment).
$ws1.Range("F4").Value = "Test LP Bingo"
$ws1.Range("G4").Value = "Houlihan Employee"
$ws1.Range("L4").Value = "United States"
$ws1.Range("M4").Value = "Active"
$ws1.Range("N4").Value = "LA"
$ws1.Range("O4").Value = "DC"
$ws1.Range("P4").Value = "Associate"
$ws1.Range("Q4").Value = "CF"
$ws1.Range("R4").Value = "CF"
$ws1.Range("S4").Value = 92001

# ---------------------------------------------------------------------------
# Sheet "ContactTypes" (sheet6): the "Archived" / "Conflicts Check LDCCR"
# rows were removed from the validation list.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("ContactTypes")
$ws6.Range("A3:A4").EntireRow.Delete()
$ws6.Range("D11").Select()

# ---------------------------------------------------------------------------
# Restore the selections on the sheets that kept theirs, and move the
# active tab from "UsersType" back to "Contact".
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A5:XFD6").Select()
